# Update countries & provincias Spain
# Refresh the COVID-19 "paises" (countries) stats table with the latest
# snapshot: bump the "last updated" timestamp, re-rank the handful of
# countries whose case counts now put them in a different position in the
# (case-count-sorted) list, and write the refreshed totals for every
# country whose figures moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 02:16"

# Re-rank countries whose updated case counts changed their sort position:
# Chequia overtakes Kenia (row 68/69)
$ws.Range("A68").Value = "Chequia"
$ws.Range("A69").Value = "Kenia"

# Gambia overtakes Somalia and Angola (rows 129-131)
$ws.Range("A129").Value = "Gambia"
$ws.Range("A130").Value = "Somalia"
$ws.Range("A131").Value = "Angola"

# Uruguay overtakes Nueva Zelanda (rows 154/155)
$ws.Range("A154").Value = "Uruguay"
$ws.Range("A155").Value = "Nueva Zelanda"

# Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
# Row 4
$ws.Range("B4").Value = 6707051
$ws.Range("C4").Value = 30450
$ws.Range("D4").Value = 3971526
$ws.Range("E4").Value = 2537049
$ws.Range("G4").Value = 348
$ws.Range("H4").Value = 198476
# Row 6
$ws.Range("E6").Value = 624834
$ws.Range("G6").Value = 389
$ws.Range("H6").Value = 131663
# Row 8
$ws.Range("B8").Value = 729619
$ws.Range("C8").Value = 6787
$ws.Range("D8").Value = 566796
$ws.Range("E8").Value = 132113
$ws.Range("G8").Value = 117
$ws.Range("H8").Value = 30710
# Row 13
$ws.Range("B13").Value = 555537
$ws.Range("C13").Value = 9056
$ws.Range("E13").Value = 124672
$ws.Range("G13").Value = 89
$ws.Range("H13").Value = 11352
# Row 36
$ws.Range("B36").Value = 101745
$ws.Range("C36").Value = 704
$ws.Range("D36").Value = 74107
$ws.Range("E36").Value = 25472
$ws.Range("G36").Value = 11
$ws.Range("H36").Value = 2166
# Row 68
$ws.Range("B68").Value = 36188
$ws.Range("C68").Value = 787
$ws.Range("D68").Value = 21294
$ws.Range("E68").Value = 14438
$ws.Range("H68").Value = 456
# Row 69
$ws.Range("B69").Value = 36157
$ws.Range("C69").Value = 188
$ws.Range("D69").Value = 23067
$ws.Range("E69").Value = 12468
$ws.Range("H69").Value = 622
# Row 74
$ws.Range("B74").Value = 27817
$ws.Range("C74").Value = 493
$ws.Range("D74").Value = 14288
$ws.Range("E74").Value = 13004
$ws.Range("G74").Value = 11
$ws.Range("H74").Value = 525
# Row 100
$ws.Range("B100").Value = 9173
$ws.Range("C100").Value = 121
$ws.Range("D100").Value = 7326
$ws.Range("E100").Value = 1815
# Row 105
$ws.Range("B105").Value = 7526
$ws.Range("C105").Value = 18
$ws.Range("D105").Value = 5678
$ws.Range("E105").Value = 1624
# Row 107
$ws.Range("B107").Value = 7238
$ws.Range("C107").Value = 44
$ws.Range("E107").Value = 717
# Row 114
$ws.Range("B114").Value = 5075
$ws.Range("C114").Value = 25
$ws.Range("D114").Value = 4229
$ws.Range("E114").Value = 748
# Row 124
$ws.Range("B124").Value = 4582
$ws.Range("C124").Value = 3
$ws.Range("D124").Value = 3850
$ws.Range("E124").Value = 639
# Row 129
$ws.Range("B129").Value = 3405
$ws.Range("C129").Value = 29
$ws.Range("D129").Value = 1723
$ws.Range("E129").Value = 1579
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 103
# Row 130
$ws.Range("B130").Value = 3389
$ws.Range("C130").Value = 13
$ws.Range("D130").Value = 2803
$ws.Range("E130").Value = 488
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 98
# Row 131
$ws.Range("B131").Value = 3388
$ws.Range("C131").Value = 53
$ws.Range("D131").Value = 1301
$ws.Range("E131").Value = 1953
$ws.Range("G131").Value = 2
$ws.Range("H131").Value = 134
# Row 154
$ws.Range("B154").Value = 1808
$ws.Range("C154").Value = 28
$ws.Range("D154").Value = 1513
$ws.Range("E154").Value = 250
$ws.Range("H154").Value = 45
# Row 155
$ws.Range("B155").Value = 1797
$ws.Range("C155").Value = 2
$ws.Range("D155").Value = 1676
$ws.Range("E155").Value = 97
$ws.Range("H155").Value = 24
# Row 158
$ws.Range("B158").Value = 1526
$ws.Range("C158").Value = 3
$ws.Range("E158").Value = 223
# Row 164
$ws.Range("B164").Value = 1180
$ws.Range("C164").Value = 2
$ws.Range("D164").Value = 1104
$ws.Range("E164").Value = 7
# Row 169
$ws.Range("D169").Value = 870
$ws.Range("E169").Value = 21
# Row 177
$ws.Range("B177").Value = 472
$ws.Range("C177").Value = 1
$ws.Range("E177").Value = 97
